$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helpers
# ---------------------------------------------------------------------------

function Find-ParaIndex($doc, $text) {
    $i = 0
    foreach ($p in $doc.Paragraphs) {
        $i += 1
        $t = $p.Range.Text.TrimEnd([char]13)
        if ($t -eq $text) {
            return $i
        }
    }
    return -1
}

# Replace the text of the (single) paragraph whose text equals $oldText with
# $newText, producing one clean run (no leftover w:proofErr / split runs).
# Works by inserting a fresh paragraph (inherits the same w:pPr) right after
# the target, filling it in, then deleting the original paragraph. Paragraph
# count is unchanged, so other indices remain valid.
function Set-ParaText($doc, $oldText, $newText) {
    $idx = Find-ParaIndex $doc $oldText
    if ($idx -eq -1) {
        throw "Set-ParaText: paragraph not found: $oldText"
    }
    $p = $doc.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $newP = $doc.Paragraphs($idx + 1)
    $newP.Range.Text = $newText
    $doc.Paragraphs($idx).Range.Delete()
}

# Insert a brand-new list paragraph with text $newText right after the
# (single) paragraph whose text equals $afterText. The new paragraph inherits
# $afterText's paragraph formatting (ListParagraph / numPr).
function Add-ParaAfter($doc, $afterText, $newText) {
    $idx = Find-ParaIndex $doc $afterText
    if ($idx -eq -1) {
        throw "Add-ParaAfter: paragraph not found: $afterText"
    }
    $p = $doc.Paragraphs($idx)
    $p.Range.InsertParagraphAfter()
    $newP = $doc.Paragraphs($idx + 1)
    $newP.Range.Text = $newText
}

# Delete the (single) paragraph whose text equals $text.
function Remove-Para($doc, $text) {
    $idx = Find-ParaIndex $doc $text
    if ($idx -eq -1) {
        throw "Remove-Para: paragraph not found: $text"
    }
    $doc.Paragraphs($idx).Range.Delete()
}

# ---------------------------------------------------------------------------
# 1) Drop the now-superfluous spell-check proofErr markers (and, where the
#    diff also shows the surrounding runs merged into one, do that too).
# ---------------------------------------------------------------------------

Set-ParaText $d "ExpressJS" "ExpressJS"
Set-ParaText $d "Formik + Yup for form and validation" "Formik + Yup for form and validation"
Set-ParaText $d "React Dropzone for image upload" "React Dropzone for image upload"
Set-ParaText $d "ExpressJS for framework" "ExpressJS for framework"
Set-ParaText $d "Multer for file upload" "Multer for file upload"

# ---------------------------------------------------------------------------
# 2) Features list edits
# ---------------------------------------------------------------------------

# Login -> Login / Logout ; Profile -> Profiles
Set-ParaText $d "Login" "Login / Logout"
Set-ParaText $d "Profile" "Profiles"

# Remove the old standalone "Search" bullet (the one between "Post Feed" and
# "Post Likes") first, before a new "Search" bullet is introduced elsewhere,
# so paragraph-text lookups stay unambiguous.
Remove-Para $d "Search"

# "Friends / Followers List" is replaced by three separate bullets:
# "Friends List", "Search", "Post Creation" (in that order, right after
# "Profiles" and before "Post Feed").
Add-ParaAfter $d "Friends / Followers List" "Friends List"
Add-ParaAfter $d "Friends List" "Search"
Add-ParaAfter $d "Search" "Post Creation"
Remove-Para $d "Friends / Followers List"

# Drop the trailing empty bullet that used to sit right after
# "Light / Dark mode".
$lightIdx = Find-ParaIndex $d "Light / Dark mode"
$trailing = $d.Paragraphs($lightIdx + 1)
if ($trailing.Range.Text.TrimEnd([char]13) -eq "") {
    $trailing.Range.Delete()
}
